$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Font/style changes -------------------------------------------------
# Title (row 1) and header row (row 2) both become bold white text
# (the header row already sits on a dark-blue fill, and now shares the
# exact same bold/white font as the title).
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215

$ws.Range("A2:K2").Font.Bold = $true
$ws.Range("A2:K2").Font.Size = 11
$ws.Range("A2:K2").Font.Color = 16777215

# --- Data changes ---------------------------------------------------------
# "PERIOD TO EXPIRE" (col H) drops by 8 (days) and "LAST UPDATE" (col I)
# moves from 08-Sep-2025 to 16-Sep-2025 for every data row (3-16).
$hValues = @{3=672;4=674;5=672;6=674;7=672;8=673;9=674;10=673;11=674;12=675;13=675;14=675;15=308;16=348}

foreach ($row in $hValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $hValues[$row]
}

# Writing the date-like text "16-Sep-2025" through .Value would make Excel
# auto-convert it to a real date serial. Use a scratch cell forced to text
# (leading apostrophe) and paste-special *values only* into each I-column
# cell so the destination keeps its original style/format while the
# content stays a literal text string, exactly like the source column.
$scratch = $ws.Cells.Item(30, 1)
$scratch.Value = "'16-Sep-2025"
$scratch.Copy()
for ($row = 3; $row -le 16; $row++) {
    $ws.Cells.Item($row, 9).PasteSpecial(-4163)
}
$scratch.Clear()
